$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.966.58'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.926.66'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.38'
$ws.Range('E6').Value = '  +5.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.930.30'
$ws.Range('E7').Value = '  +2.24%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('E11').Value = '  +1.74%  '
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('E13').Value = '  +5.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.39'
$ws.Range('E14').Value = '  +4.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.585.47'
$ws.Range('E15').Value = '  +2.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.939.32'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.941.51'
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.71'
$ws.Range('E18').Value = '  +9.36%  '
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.17'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.63'
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.748'
$ws.Range('E23').Value = '  +4.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000168'
$ws.Range('E24').Value = '  +5.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.95'
$ws.Range('E25').Value = '  +2.44%  '
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.33'
$ws.Range('E27').Value = '  +2.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.13'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('E31').Value = '  +3.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.078.77'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '32.19'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.890.55'
$ws.Range('E35').Value = '  +2.59%  '
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('E37').Value = '  +4.27%  '
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('E40').Value = '  +10.83%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.330'
$ws.Range('E41').Value = '  +3.52%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.13'
$ws.Range('E43').Value = '  +7.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '436.50'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.66'
$ws.Range('E46').Value = '  +3.26%  '
$ws.Range('E48').Value = '  +3.62%  '
$ws.Range('E49').Value = '  +23.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '40.57'
$ws.Range('E50').Value = '  +5.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.07'
$ws.Range('E51').Value = '  -0.17%  '
